$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2022-05-30 T 18:36:33 UTC"
$ws.Range("B2").Value = 30619.07644016125
$ws.Range("C2").Value = 1.00005
$ws.Range("D2").Value = 1.226498

$ws.Range("A3").Value = "2022-05-30 T 18:36:33 UTC"
$ws.Range("B3").Value = 30619.07644016125
$ws.Range("C3").Value = 1.00005
$ws.Range("D3").Value = 1.226498
